# Release refresh: the embedded build timestamp in the version strings
# moved from "February 03 2026 17.29.55 EST" to "February 03 2026 18.05.36 EST".
# Update every cell that carries that timestamp: the "About" sheet's
# version banner (A2) and citation text (A6), plus the per-row version
# column (S2:S11) on the "Boundaries and methane sources" sheet.

$wb = $excel.ActiveWorkbook

$oldTime = "17.29.55"
$newTime = "18.05.36"

$aboutSheet = $wb.Worksheets.Item("About")
$dataSheet  = $wb.Worksheets.Item("Boundaries and methane sources")

$targets = @($aboutSheet.Range("A2"), $aboutSheet.Range("A6"))
for ($row = 2; $row -le 11; $row++) {
    $targets += $dataSheet.Range("S$row")
}

foreach ($cell in $targets) {
    $text = $cell.Value2
    if ($text -ne $null -and $text -is [string] -and $text.Contains($oldTime)) {
        $cell.Value2 = $text.Replace($oldTime, $newTime)
    }
}
